# "Allow external senders for position emails"
#
# The schema sheet documents each DynamoDB-style item shape as a PK/SK block.
# This change removes the now-unused "user_email_2" alias block, and reuses
# the vacated rows for the "list_email_1" block (which gains a new
# "allow_external" attribute), plus a brand-new "role_id" Other-Attributes
# row (role -> description). The now-orphaned "list_email_2" block is
# dropped in the process.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Delete the "user_email_2" entity block (rows 7-9, merged B7:B9).
# Excel's row delete shifts everything below up, so the "list_email_1"
# block (old rows 10-12) becomes rows 7-9, and the lone "list_email_2"
# row (old row 13) becomes row 10 - this also drops the orphaned
# "user_email_2" / "list_email_2" strings and the B10:B12 merge.
$ws.Rows("7:9").Delete()

# Step 2: Turn the old "list_email_2" row (now row 10) into the new
# "role_id" Other-Attributes row: role_id references "role" and carries
# a "description".
$ws.Range("B10").Value = "role_id"
$ws.Range("C10").Value = "role"
$ws.Range("D10").Value = "description"

# Step 3: Add the new "allow_external" attribute to the "list_email_1"
# block header (now row 7), letting positions opt in to accepting mail
# from senders outside the list.
$ws.Range("G7").Value = "allow_external"

# Step 4: Match the author's final selection.
$ws.Range("G8").Select() | Out-Null
